# Trade #21 closed at 2026-02-17 15:18:51 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.77
$wsSummary.Range("B4").Value = -0.24
$wsSummary.Range("B6").Value = 21
$wsSummary.Range("B8").Value = 11
$wsSummary.Range("B9").Value = 23.81

# ---- Strategy Status sheet (MarketMaking row, row 4) ----
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.77
$wsStatus.Range("D4").Value = 21
$wsStatus.Range("E4").Value = -0.24
$wsStatus.Range("F4").Value = -0.23
$wsStatus.Range("G4").Value = 23.81

# ---- Add new closed trade (row 22) to "All Trades" and "MarketMaking" sheets ----
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(22, 1).Value = 21

    # Column B holds a date-like string ("2026-02-17"); force text formatting
    # so Excel keeps it as literal text instead of auto-converting it to a
    # date serial number, then drop back to the Normal style so no extra
    # per-cell formatting is left behind.
    $ws.Cells.Item(22, 2).NumberFormat = "@"
    $ws.Cells.Item(22, 2).Value = "2026-02-17"
    $ws.Cells.Item(22, 2).Style = "Normal"

    $ws.Cells.Item(22, 3).Value = "15:18:45"
    $ws.Cells.Item(22, 4).Value = "MarketMaking"
    $ws.Cells.Item(22, 5).Value = "DOWN"
    $ws.Cells.Item(22, 6).Value = 0.18
    $ws.Cells.Item(22, 7).Value = 0.17
    $ws.Cells.Item(22, 8).Value = "CLOSED"
    $ws.Cells.Item(22, 9).Value = -5.5556
    $ws.Cells.Item(22, 10).Value = -0.01
    $ws.Cells.Item(22, 11).Value = 99.77
    $ws.Cells.Item(22, 12).Value = 0
    $ws.Cells.Item(22, 13).Value = 0
    $ws.Cells.Item(22, 14).Value = 0.6
    $ws.Cells.Item(22, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(22, 16).Value = "early_exit"
    $ws.Cells.Item(22, 17).Value = 0.14
}
